# Rename column headers from "_old"/"_new" suffixes to "_FV2410"/"_FV2504"
# and turn the used range into a real Excel Table ("Table1") with autofilter,
# then freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into a proper ListObject (Excel Table) with an
# autofilter, matching the original data extent A1:U88.
$tableRange = $ws.Range("A1:U88")
$table = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$table.Name = "Table1"

# Freeze the header row (split after row 1, keep column A visible).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
